$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("financial_statements")

# Step 1: Insert new rows to make room for new dictionary entries.
# Insert 1 row at row 24 (new: accumulated_depreciation_and_amortization)
$ws.Range("A24:A24").EntireRow.Insert()
# Insert 6 rows at row 54 (new equity breakdown fields)
$ws.Range("A54:A59").EntireRow.Insert()
# Insert 2 rows at row 78 (new income statement financial gains/costs fields)
$ws.Range("A78:A79").EntireRow.Insert()

# Step 2: Write full row contents (A:I) for rows 2-88 with final values.
$data = New-Object 'object[,]' 87,9
$data[0,0] = 'id'
$data[0,1] = 'Identificador único de Belvo para el elemento actual.'
$data[0,2] = '0d3ffb69-f83b-456e-ad8e-208d0998d71d'
$data[0,3] = 'string'
$data[0,4] = 'uuid'
$data[0,5] = 'Yes'
$data[0,6] = 'No'
$data[0,7] = $null
$data[0,8] = $null
$data[1,0] = 'link'
$data[1,1] = 'El `link.id` al que pertenecen los datos.'
$data[1,2] = '30cb4806-6e00-48a4-91c9-ca55968576c8'
$data[1,3] = 'string'
$data[1,4] = 'uuid'
$data[1,5] = 'Yes'
$data[1,6] = 'Yes'
$data[1,7] = $null
$data[1,8] = $null
$data[2,0] = 'collected_at'
$data[2,1] = 'La marca de tiempo ISO-8601 cuando se recopiló el punto de datos.'
$data[2,2] = '2022-02-09T08:45:50.406032Z'
$data[2,3] = 'string'
$data[2,4] = 'date-time'
$data[2,5] = 'Yes'
$data[2,6] = 'No'
$data[2,7] = $null
$data[2,8] = $null
$data[3,0] = 'created_at'
$data[3,1] = 'La marca de tiempo ISO-8601 de cuando se creó el punto de datos en la base de datos de Belvo.'
$data[3,2] = '2022-02-09T08:45:50.406032Z'
$data[3,3] = 'string'
$data[3,4] = 'date-time'
$data[3,5] = 'Yes'
$data[3,6] = 'No'
$data[3,7] = $null
$data[3,8] = $null
$data[4,0] = 'error'
$data[4,1] = 'En casos donde surjan problemas durante la extracción de estados financieros de la institución fiscal, se pueden proporcionar los siguientes mensajes de error para explicar los problemas encontrados:

  - `Unable to validate if the user has an available financial statement for the specified year.`
  - `No available financial statement found for the user for the specified year, preventing data extraction.`
  - `Unable to verify if the user has *conceptos vigentes* for the specified year.`
  - `The fiscal institution provided the financial statement in an unrecognized format.`'
$data[4,2] = $null
$data[4,3] = 'string'
$data[4,4] = $null
$data[4,5] = 'Yes'
$data[4,6] = 'Yes'
$data[4,7] = $null
$data[4,8] = $null
$data[5,0] = 'year'
$data[5,1] = 'El año del estado financiero.'
$data[5,2] = '2020'
$data[5,3] = 'string'
$data[5,4] = $null
$data[5,5] = 'Yes'
$data[5,6] = 'No'
$data[5,7] = $null
$data[5,8] = '^\d{4}$'
$data[6,0] = 'currency'
$data[6,1] = 'La moneda del estado financiero.'
$data[6,2] = 'MXN'
$data[6,3] = 'string'
$data[6,4] = $null
$data[6,5] = 'Yes'
$data[6,6] = 'No'
$data[6,7] = $null
$data[6,8] = $null
$data[7,0] = 'balance_sheet'
$data[7,1] = 'El balance general que detalla los activos, pasivos y patrimonio de la empresa para el año dado.'
$data[7,2] = $null
$data[7,3] = 'object'
$data[7,4] = $null
$data[7,5] = 'Yes'
$data[7,6] = 'No'
$data[7,7] = $null
$data[7,8] = $null
$data[8,0] = 'balance_sheet.current_assets'
$data[8,1] = 'Los activos corrientes de la empresa para el año dado.'
$data[8,2] = $null
$data[8,3] = 'object'
$data[8,4] = $null
$data[8,5] = 'No'
$data[8,6] = 'No'
$data[8,7] = $null
$data[8,8] = $null
$data[9,0] = 'balance_sheet.current_assets.cash_and_equivalents'
$data[9,1] = 'El monto total de efectivo y equivalentes de efectivo, incluyendo moneda, cuentas bancarias y otras inversiones líquidas que pueden convertirse rápidamente en efectivo.'
$data[9,2] = '48572.01'
$data[9,3] = 'number'
$data[9,4] = 'float'
$data[9,5] = 'Yes'
$data[9,6] = 'Yes'
$data[9,7] = $null
$data[9,8] = $null
$data[10,0] = 'balance_sheet.current_assets.short_term_investments'
$data[10,1] = 'El valor de las inversiones que se espera sean liquidadas en efectivo dentro de un año, como los valores negociables.'
$data[10,2] = '21345.01'
$data[10,3] = 'number'
$data[10,4] = 'float'
$data[10,5] = 'Yes'
$data[10,6] = 'Yes'
$data[10,7] = $null
$data[10,8] = $null
$data[11,0] = 'balance_sheet.current_assets.accounts_receivable'
$data[11,1] = 'El monto adeudado por los clientes por ventas realizadas a crédito, que se espera recibir en un corto período.'
$data[11,2] = '154321.01'
$data[11,3] = 'number'
$data[11,4] = 'float'
$data[11,5] = 'Yes'
$data[11,6] = 'Yes'
$data[11,7] = $null
$data[11,8] = $null
$data[12,0] = 'balance_sheet.current_assets.notes_receivable'
$data[12,1] = 'El valor de los pagarés escritos recibidos de clientes u otros, que prometen pagar una cantidad especificada en una fecha determinada.'
$data[12,2] = '31789.01'
$data[12,3] = 'number'
$data[12,4] = 'float'
$data[12,5] = 'Yes'
$data[12,6] = 'Yes'
$data[12,7] = $null
$data[12,8] = $null
$data[13,0] = 'balance_sheet.current_assets.other_debtors'
$data[13,1] = 'Los montos totales adeudados por varios otros deudores, excluyendo cuentas y documentos por cobrar.'
$data[13,2] = '12345.01'
$data[13,3] = 'number'
$data[13,4] = 'float'
$data[13,5] = 'Yes'
$data[13,6] = 'Yes'
$data[13,7] = $null
$data[13,8] = $null
$data[14,0] = 'balance_sheet.current_assets.bad_debt_provision'
$data[14,1] = 'El monto estimado de cuentas por cobrar que se espera que sean incobrables, a menudo se denomina provisión para cuentas de dudoso cobro.'
$data[14,2] = '0.01'
$data[14,3] = 'number'
$data[14,4] = 'float'
$data[14,5] = 'Yes'
$data[14,6] = 'Yes'
$data[14,7] = $null
$data[14,8] = $null
$data[15,0] = 'balance_sheet.current_assets.tax_recoverable'
$data[15,1] = 'El monto de los pagos de impuestos que se pueden recuperar de las autoridades fiscales.'
$data[15,2] = '8976.01'
$data[15,3] = 'number'
$data[15,4] = 'float'
$data[15,5] = 'Yes'
$data[15,6] = 'Yes'
$data[15,7] = $null
$data[15,8] = $null
$data[16,0] = 'balance_sheet.current_assets.inventory'
$data[16,1] = 'El valor total de los bienes disponibles para la venta, materias primas, trabajo en proceso y productos terminados.'
$data[16,2] = '65432.01'
$data[16,3] = 'number'
$data[16,4] = 'float'
$data[16,5] = 'Yes'
$data[16,6] = 'Yes'
$data[16,7] = $null
$data[16,8] = $null
$data[17,0] = 'balance_sheet.current_assets.prepaid_expenses'
$data[17,1] = 'El monto pagado por adelantado por bienes o servicios que se recibirán en el futuro, como primas de seguro o alquiler.'
$data[17,2] = '14321.01'
$data[17,3] = 'number'
$data[17,4] = 'float'
$data[17,5] = 'Yes'
$data[17,6] = 'Yes'
$data[17,7] = $null
$data[17,8] = $null
$data[18,0] = 'balance_sheet.current_assets.assets_available_for_sale'
$data[18,1] = 'El valor de los activos no corrientes que están disponibles para la venta pero aún no se han vendido, como equipos excedentes o propiedades.'
$data[18,2] = '54321.01'
$data[18,3] = 'number'
$data[18,4] = 'float'
$data[18,5] = 'Yes'
$data[18,6] = 'Yes'
$data[18,7] = $null
$data[18,8] = $null
$data[19,0] = 'balance_sheet.current_assets.total'
$data[19,1] = 'La suma de todos los activos corrientes, que representa el valor total de los activos que se espera convertir en efectivo o utilizar dentro de un año.'
$data[19,2] = '372480.01'
$data[19,3] = 'number'
$data[19,4] = 'float'
$data[19,5] = 'Yes'
$data[19,6] = 'Yes'
$data[19,7] = $null
$data[19,8] = $null
$data[20,0] = 'balance_sheet.non_current_assets'
$data[20,1] = 'Los activos no corrientes de la empresa, que son inversiones a largo plazo o propiedades que no se convierten fácilmente en efectivo, para el año dado.'
$data[20,2] = $null
$data[20,3] = 'object'
$data[20,4] = $null
$data[20,5] = 'No'
$data[20,6] = 'No'
$data[20,7] = $null
$data[20,8] = $null
$data[21,0] = 'balance_sheet.non_current_assets.property_plant_and_equipment'
$data[21,1] = 'El valor total de los bienes inmuebles, planta y equipo propiedad de la empresa, incluyendo terrenos, edificios, maquinaria y vehículos, utilizados para operaciones a largo plazo.'
$data[21,2] = '1123456.01'
$data[21,3] = 'number'
$data[21,4] = 'float'
$data[21,5] = 'Yes'
$data[21,6] = 'Yes'
$data[21,7] = $null
$data[21,8] = $null
$data[22,0] = 'balance_sheet.non_current_assets.accumulated_depreciation_and_amortization'
$data[22,1] = 'La depreciación y amortización acumulada total, que representa la asignación acumulativa del costo de los activos no corrientes durante el período en que se espera que proporcionen beneficios económicos.'
$data[22,2] = '123456.01'
$data[22,3] = 'number'
$data[22,4] = 'float'
$data[22,5] = 'Yes'
$data[22,6] = 'Yes'
$data[22,7] = $null
$data[22,8] = $null
$data[23,0] = 'balance_sheet.non_current_assets.long_term_accounts_receivable'
$data[23,1] = 'El monto adeudado por los clientes por ventas realizadas a crédito, que se espera recibir después de un año.'
$data[23,2] = '10987.01'
$data[23,3] = 'number'
$data[23,4] = 'float'
$data[23,5] = 'Yes'
$data[23,6] = 'Yes'
$data[23,7] = $null
$data[23,8] = $null
$data[24,0] = 'balance_sheet.non_current_assets.prepayment_to_suppliers'
$data[24,1] = 'El monto pagado por adelantado a los proveedores por bienes o servicios que se recibirán en el futuro, que se espera utilizar a largo plazo.'
$data[24,2] = '5432.01'
$data[24,3] = 'number'
$data[24,4] = 'float'
$data[24,5] = 'Yes'
$data[24,6] = 'Yes'
$data[24,7] = $null
$data[24,8] = $null
$data[25,0] = 'balance_sheet.non_current_assets.goodwill'
$data[25,1] = 'El valor de los activos intangibles que surgen de la adquisición de otras empresas, representando la prima pagada sobre el valor razonable de los activos netos adquiridos.'
$data[25,2] = '47654.01'
$data[25,3] = 'number'
$data[25,4] = 'float'
$data[25,5] = 'Yes'
$data[25,6] = 'Yes'
$data[25,7] = $null
$data[25,8] = $null
$data[26,0] = 'balance_sheet.non_current_assets.intangible_assets'
$data[26,1] = 'El valor total de los activos intangibles propiedad de la empresa, como patentes, marcas registradas y derechos de autor, con vidas útiles que se extienden más allá de un año.'
$data[26,2] = '43210.01'
$data[26,3] = 'number'
$data[26,4] = 'float'
$data[26,5] = 'Yes'
$data[26,6] = 'Yes'
$data[26,7] = $null
$data[26,8] = $null
$data[27,0] = 'balance_sheet.non_current_assets.investments_in_associates'
$data[27,1] = 'El valor de las inversiones en otras empresas en las que la compañía tiene una influencia significativa pero no control, típicamente representado por la propiedad del 20-50% de las acciones con derecho a voto del asociado.'
$data[27,2] = '65432.01'
$data[27,3] = 'number'
$data[27,4] = 'float'
$data[27,5] = 'Yes'
$data[27,6] = 'Yes'
$data[27,7] = $null
$data[27,8] = $null
$data[28,0] = 'balance_sheet.non_current_assets.long_term_financial_instruments'
$data[28,1] = 'El valor de los instrumentos financieros que se espera mantener por más de un año, como bonos, debentures y préstamos a largo plazo.'
$data[28,2] = '32876.01'
$data[28,3] = 'number'
$data[28,4] = 'float'
$data[28,5] = 'Yes'
$data[28,6] = 'Yes'
$data[28,7] = $null
$data[28,8] = $null
$data[29,0] = 'balance_sheet.non_current_assets.total'
$data[29,1] = 'La suma de todos los activos no corrientes, que representa el valor total de los activos que se espera utilizar o mantener por más de un año.'
$data[29,2] = '1346647.01'
$data[29,3] = 'number'
$data[29,4] = 'float'
$data[29,5] = 'Yes'
$data[29,6] = 'Yes'
$data[29,7] = $null
$data[29,8] = $null
$data[30,0] = 'balance_sheet.current_liabilities'
$data[30,1] = 'Los pasivos corrientes de la empresa, que se espera sean liquidados dentro del año dado.'
$data[30,2] = $null
$data[30,3] = 'object'
$data[30,4] = $null
$data[30,5] = 'No'
$data[30,6] = 'No'
$data[30,7] = $null
$data[30,8] = $null
$data[31,0] = 'balance_sheet.current_liabilities.bank_loans'
$data[31,1] = 'El monto total de los préstamos tomados de bancos o instituciones financieras, que se espera sean reembolsados dentro de un año.'
$data[31,2] = '49876.01'
$data[31,3] = 'number'
$data[31,4] = 'float'
$data[31,5] = 'Yes'
$data[31,6] = 'Yes'
$data[31,7] = $null
$data[31,8] = $null
$data[32,0] = 'balance_sheet.current_liabilities.accounts_payable'
$data[32,1] = 'El monto adeudado a los proveedores por bienes o servicios comprados a crédito, que se espera pagar en un corto período.'
$data[32,2] = '103298.01'
$data[32,3] = 'number'
$data[32,4] = 'float'
$data[32,5] = 'Yes'
$data[32,6] = 'Yes'
$data[32,7] = $null
$data[32,8] = $null
$data[33,0] = 'balance_sheet.current_liabilities.notes_payable'
$data[33,1] = 'El valor de los pagarés emitidos a proveedores u otros, prometiendo pagar una cantidad especificada en una fecha determinada.'
$data[33,2] = '25643.01'
$data[33,3] = 'number'
$data[33,4] = 'float'
$data[33,5] = 'Yes'
$data[33,6] = 'Yes'
$data[33,7] = $null
$data[33,8] = $null
$data[34,0] = 'balance_sheet.current_liabilities.financial_instruments'
$data[34,1] = 'El valor de los instrumentos financieros que se espera sean liquidados en efectivo dentro de un año, como bonos, debentures y préstamos a corto plazo.'
$data[34,2] = '14321.01'
$data[34,3] = 'number'
$data[34,4] = 'float'
$data[34,5] = 'Yes'
$data[34,6] = 'Yes'
$data[34,7] = $null
$data[34,8] = $null
$data[35,0] = 'balance_sheet.current_liabilities.other_creditors'
$data[35,1] = 'Los montos totales adeudados a varios otros acreedores, excluyendo cuentas y pagarés por pagar.'
$data[35,2] = '21987.01'
$data[35,3] = 'number'
$data[35,4] = 'float'
$data[35,5] = 'Yes'
$data[35,6] = 'Yes'
$data[35,7] = $null
$data[35,8] = $null
$data[36,0] = 'balance_sheet.current_liabilities.income_tax_payable'
$data[36,1] = 'La cantidad de impuesto sobre la renta que se debe a las autoridades fiscales, que se espera pagar en un corto período.'
$data[36,2] = '12765.01'
$data[36,3] = 'number'
$data[36,4] = 'float'
$data[36,5] = 'Yes'
$data[36,6] = 'Yes'
$data[36,7] = $null
$data[36,8] = $null
$data[37,0] = 'balance_sheet.current_liabilities.customer_advances'
$data[37,1] = 'El monto total recibido por adelantado de los clientes por bienes o servicios que se entregarán en el futuro, que se espera utilizar dentro de un año.'
$data[37,2] = '18765.01'
$data[37,3] = 'number'
$data[37,4] = 'float'
$data[37,5] = 'Yes'
$data[37,6] = 'Yes'
$data[37,7] = $null
$data[37,8] = $null
$data[38,0] = 'balance_sheet.current_liabilities.provisions'
$data[38,1] = 'El monto estimado reservado para pasivos o pérdidas futuras, como garantías, reclamaciones legales o costos de reestructuración.'
$data[38,2] = '10987.01'
$data[38,3] = 'number'
$data[38,4] = 'float'
$data[38,5] = 'Yes'
$data[38,6] = 'Yes'
$data[38,7] = $null
$data[38,8] = $null
$data[39,0] = 'balance_sheet.current_liabilities.taxes_payable'
$data[39,1] = 'El monto total de impuestos adeudados a las autoridades fiscales, que se espera pagar en un corto período.'
$data[39,2] = '5321.01'
$data[39,3] = 'number'
$data[39,4] = 'float'
$data[39,5] = 'Yes'
$data[39,6] = 'Yes'
$data[39,7] = $null
$data[39,8] = $null
$data[40,0] = 'balance_sheet.current_liabilities.total'
$data[40,1] = 'La suma de todos los pasivos corrientes, que representa el valor total de las obligaciones que se espera liquidar dentro de un año.'
$data[40,2] = '260963.01'
$data[40,3] = 'number'
$data[40,4] = 'float'
$data[40,5] = 'Yes'
$data[40,6] = 'Yes'
$data[40,7] = $null
$data[40,8] = $null
$data[41,0] = 'balance_sheet.non_current_liabilities'
$data[41,1] = 'Los pasivos no corrientes de la empresa, que son obligaciones a largo plazo que no vencen dentro del año dado.'
$data[41,2] = $null
$data[41,3] = 'object'
$data[41,4] = $null
$data[41,5] = 'No'
$data[41,6] = 'No'
$data[41,7] = $null
$data[41,8] = $null
$data[42,0] = 'balance_sheet.non_current_liabilities.long_term_accounts_payable'
$data[42,1] = 'El monto adeudado a los proveedores por bienes o servicios comprados a crédito, que se espera pagar después de un año.'
$data[42,2] = '30876.01'
$data[42,3] = 'number'
$data[42,4] = 'float'
$data[42,5] = 'Yes'
$data[42,6] = 'Yes'
$data[42,7] = $null
$data[42,8] = $null
$data[43,0] = 'balance_sheet.non_current_liabilities.long_term_financial_instruments'
$data[43,1] = 'El valor de los instrumentos financieros que se espera mantener por más de un año, como bonos, debentures y préstamos a largo plazo.'
$data[43,2] = '42310.01'
$data[43,3] = 'number'
$data[43,4] = 'float'
$data[43,5] = 'Yes'
$data[43,6] = 'Yes'
$data[43,7] = $null
$data[43,8] = $null
$data[44,0] = 'balance_sheet.non_current_liabilities.deferred_revenue'
$data[44,1] = 'El monto recibido por adelantado de los clientes por bienes o servicios que se entregarán en el futuro, que se espera reconocer como ingresos a largo plazo (como el alquiler).'
$data[44,2] = '21987.01'
$data[44,3] = 'number'
$data[44,4] = 'float'
$data[44,5] = 'Yes'
$data[44,6] = 'Yes'
$data[44,7] = $null
$data[44,8] = $null
$data[45,0] = 'balance_sheet.non_current_liabilities.contributions_for_future_capital_increases'
$data[45,1] = 'Las contribuciones totales recibidas de los accionistas u otros inversores para futuros aumentos de capital, que se espera sean utilizadas a largo plazo.'
$data[45,2] = '10987.01'
$data[45,3] = 'number'
$data[45,4] = 'float'
$data[45,5] = 'Yes'
$data[45,6] = 'Yes'
$data[45,7] = $null
$data[45,8] = $null
$data[46,0] = 'balance_sheet.non_current_liabilities.deferred_income_tax'
$data[46,1] = 'La cantidad de impuesto sobre la renta que se difiere a períodos futuros, que se espera pagar después de un año.'
$data[46,2] = '26543.01'
$data[46,3] = 'number'
$data[46,4] = 'float'
$data[46,5] = 'Yes'
$data[46,6] = 'Yes'
$data[46,7] = $null
$data[46,8] = $null
$data[47,0] = 'balance_sheet.non_current_liabilities.employee_benefits'
$data[47,1] = 'El monto total de los beneficios adeudados a los empleados, como pensiones, gratificaciones y otros beneficios posteriores al empleo, que se espera liquidar a largo plazo.'
$data[47,2] = '30218.01'
$data[47,3] = 'number'
$data[47,4] = 'float'
$data[47,5] = 'Yes'
$data[47,6] = 'Yes'
$data[47,7] = $null
$data[47,8] = $null
$data[48,0] = 'balance_sheet.non_current_liabilities.long_term_provisions'
$data[48,1] = 'El monto estimado reservado para pasivos o pérdidas futuras, como garantías, reclamaciones legales o costos de reestructuración, que se espera liquidar después de un año.'
$data[48,2] = '15432.01'
$data[48,3] = 'number'
$data[48,4] = 'float'
$data[48,5] = 'Yes'
$data[48,6] = 'Yes'
$data[48,7] = $null
$data[48,8] = $null
$data[49,0] = 'balance_sheet.non_current_liabilities.total'
$data[49,1] = 'La suma de todos los pasivos no corrientes, que representa el valor total de las obligaciones que se espera liquidar después de un año.'
$data[49,2] = '178353.01'
$data[49,3] = 'number'
$data[49,4] = 'float'
$data[49,5] = 'Yes'
$data[49,6] = 'Yes'
$data[49,7] = $null
$data[49,8] = $null
$data[50,0] = 'balance_sheet.equity'
$data[50,1] = 'El patrimonio de la empresa, que representa el interés residual en los activos después de deducir los pasivos.'
$data[50,2] = $null
$data[50,3] = 'object'
$data[50,4] = $null
$data[50,5] = 'No'
$data[50,6] = 'No'
$data[50,7] = $null
$data[50,8] = $null
$data[51,0] = 'balance_sheet.equity.stockholders_equity'
$data[51,1] = 'El valor total de las acciones emitidas por la empresa, que representa el interés de propiedad de los accionistas en el negocio.'
$data[51,2] = '501234.01'
$data[51,3] = 'number'
$data[51,4] = 'float'
$data[51,5] = 'Yes'
$data[51,6] = 'Yes'
$data[51,7] = $null
$data[51,8] = $null
$data[52,0] = 'balance_sheet.equity.future_capital_contributions'
$data[52,1] = 'Los fondos recibidos de los accionistas que están específicamente designados para futuros aumentos de capital o inversiones.'
$data[52,2] = '75000.01'
$data[52,3] = 'number'
$data[52,4] = 'float'
$data[52,5] = 'Yes'
$data[52,6] = 'Yes'
$data[52,7] = $null
$data[52,8] = $null
$data[53,0] = 'balance_sheet.equity.legal_reserve'
$data[53,1] = 'La reserva legal exigida por la ley, generalmente apartada de las ganancias, para proporcionar protección financiera contra pérdidas u obligaciones futuras.'
$data[53,2] = '25000.01'
$data[53,3] = 'number'
$data[53,4] = 'float'
$data[53,5] = 'Yes'
$data[53,6] = 'Yes'
$data[53,7] = $null
$data[53,8] = $null
$data[54,0] = 'balance_sheet.equity.capital_update_excess'
$data[54,1] = 'El excedente resultante de los ajustes realizados al capital social, a menudo debido a la inflación o la revalorización de activos.'
$data[54,2] = '15000.01'
$data[54,3] = 'number'
$data[54,4] = 'float'
$data[54,5] = 'Yes'
$data[54,6] = 'Yes'
$data[54,7] = $null
$data[54,8] = $null
$data[55,0] = 'balance_sheet.equity.capital_update_insufficiency'
$data[55,1] = 'El déficit resultante de los ajustes realizados al capital social, a menudo debido a la inflación o la revalorización de activos.'
$data[55,2] = '-5000.01'
$data[55,3] = 'number'
$data[55,4] = 'float'
$data[55,5] = 'Yes'
$data[55,6] = 'Yes'
$data[55,7] = $null
$data[55,8] = $null
$data[56,0] = 'balance_sheet.equity.capital_reserve'
$data[56,1] = 'La reserva de capital derivada de actividades no operativas, como ganancias de revalorizaciones de activos o ciertas transacciones de capital.'
$data[56,2] = '10000.01'
$data[56,3] = 'number'
$data[56,4] = 'float'
$data[56,5] = 'Yes'
$data[56,6] = 'Yes'
$data[56,7] = $null
$data[56,8] = $null
$data[57,0] = 'balance_sheet.equity.share_premium_on_stock_sales'
$data[57,1] = 'El monto excedente recibido por una empresa cuando las acciones se emiten a un precio superior a su valor nominal (par).'
$data[57,2] = '50000.01'
$data[57,3] = 'number'
$data[57,4] = 'float'
$data[57,5] = 'Yes'
$data[57,6] = 'Yes'
$data[57,7] = $null
$data[57,8] = $null
$data[58,0] = 'balance_sheet.equity.retained_earnings'
$data[58,1] = 'Las ganancias o pérdidas acumuladas de la empresa que no se han distribuido a los accionistas como dividendos.'
$data[58,2] = '202345.01'
$data[58,3] = 'number'
$data[58,4] = 'float'
$data[58,5] = 'Yes'
$data[58,6] = 'Yes'
$data[58,7] = $null
$data[58,8] = $null
$data[59,0] = 'balance_sheet.equity.other_comprehensive_income'
$data[59,1] = 'Las ganancias o pérdidas que no se incluyen en el ingreso neto pero se reportan directamente en el patrimonio, como las ganancias no realizadas en inversiones o los ajustes por conversión de moneda extranjera.'
$data[59,2] = '10987.01'
$data[59,3] = 'number'
$data[59,4] = 'float'
$data[59,5] = 'Yes'
$data[59,6] = 'Yes'
$data[59,7] = $null
$data[59,8] = $null
$data[60,0] = 'balance_sheet.equity.controlling_interest'
$data[60,1] = 'El interés de propiedad en la empresa que posee la entidad matriz o los accionistas mayoritarios, que representa la participación de control en el negocio.'
$data[60,2] = '70876.01'
$data[60,3] = 'number'
$data[60,4] = 'float'
$data[60,5] = 'Yes'
$data[60,6] = 'Yes'
$data[60,7] = $null
$data[60,8] = $null
$data[61,0] = 'balance_sheet.equity.non_controlling_interest'
$data[61,1] = 'La participación accionaria en la empresa mantenida por los accionistas minoritarios, que representa la participación no controladora en el negocio.'
$data[61,2] = '50321.01'
$data[61,3] = 'number'
$data[61,4] = 'float'
$data[61,5] = 'Yes'
$data[61,6] = 'Yes'
$data[61,7] = $null
$data[61,8] = $null
$data[62,0] = 'balance_sheet.equity.total'
$data[62,1] = 'La suma del capital social, las ganancias retenidas, otros ingresos integrales, el interés controlador y el interés no controlador, que representa el patrimonio total de la empresa.'
$data[62,2] = '836763.01'
$data[62,3] = 'number'
$data[62,4] = 'float'
$data[62,5] = 'Yes'
$data[62,6] = 'Yes'
$data[62,7] = $null
$data[62,8] = $null
$data[63,0] = 'income_statement'
$data[63,1] = 'El estado de resultados que detalla los ingresos, gastos y beneficios de la empresa para el año dado.'
$data[63,2] = $null
$data[63,3] = 'object'
$data[63,4] = $null
$data[63,5] = 'Yes'
$data[63,6] = 'No'
$data[63,7] = $null
$data[63,8] = $null
$data[64,0] = 'income_statement.net_revenue'
$data[64,1] = 'Los ingresos totales generados por la empresa a partir de sus operaciones comerciales principales, excluyendo cualquier deducción por descuentos, devoluciones o bonificaciones.

> **Nota**: `domestic_sales` + `foreign_sales` no sumarán el `net_revenue` debido a la exclusión de descuentos, devoluciones y bonificaciones.'
$data[64,2] = '1212345.01'
$data[64,3] = 'number'
$data[64,4] = 'float'
$data[64,5] = 'Yes'
$data[64,6] = 'Yes'
$data[64,7] = $null
$data[64,8] = $null
$data[65,0] = 'income_statement.domestic_sales'
$data[65,1] = 'Los ingresos generados por la empresa a partir de la venta de bienes o servicios dentro de su país de origen.'
$data[65,2] = '1123456.01'
$data[65,3] = 'number'
$data[65,4] = 'float'
$data[65,5] = 'Yes'
$data[65,6] = 'Yes'
$data[65,7] = $null
$data[65,8] = $null
$data[66,0] = 'income_statement.foreign_sales'
$data[66,1] = 'Los ingresos generados por la empresa a partir de la venta de bienes o servicios en países extranjeros.'
$data[66,2] = '88987.01'
$data[66,3] = 'number'
$data[66,4] = 'float'
$data[66,5] = 'Yes'
$data[66,6] = 'Yes'
$data[66,7] = $null
$data[66,8] = $null
$data[67,0] = 'income_statement.materials_used'
$data[67,1] = 'El costo total de los materiales utilizados o comercializados por la empresa durante el período de informe.'
$data[67,2] = '609876.01'
$data[67,3] = 'number'
$data[67,4] = 'float'
$data[67,5] = 'No'
$data[67,6] = 'Yes'
$data[67,7] = $null
$data[67,8] = $null
$data[68,0] = 'income_statement.cost_of_goods_sold'
$data[68,1] = 'El costo total incurrido por la empresa para producir o comprar los bienes vendidos durante el período de informe.'
$data[68,2] = '412345.01'
$data[68,3] = 'number'
$data[68,4] = 'float'
$data[68,5] = 'Yes'
$data[68,6] = 'Yes'
$data[68,7] = $null
$data[68,8] = $null
$data[69,0] = 'income_statement.cost_of_services_sold'
$data[69,1] = 'El costo total incurrido por la empresa para proporcionar los servicios vendidos durante el período de informe.'
$data[69,2] = '101234.01'
$data[69,3] = 'number'
$data[69,4] = 'float'
$data[69,5] = 'Yes'
$data[69,6] = 'Yes'
$data[69,7] = $null
$data[69,8] = $null
$data[70,0] = 'income_statement.gross_profit'
$data[70,1] = 'La diferencia entre los ingresos netos y el costo total de los bienes y servicios vendidos, que representa la ganancia obtenida de las operaciones comerciales principales antes de deducir los gastos operativos.'
$data[70,2] = '190890.01'
$data[70,3] = 'number'
$data[70,4] = 'float'
$data[70,5] = 'Yes'
$data[70,6] = 'Yes'
$data[70,7] = $null
$data[70,8] = $null
$data[71,0] = 'income_statement.gross_loss'
$data[71,1] = 'La diferencia negativa entre los ingresos netos y el costo total de bienes y servicios vendidos, que representa la pérdida incurrida por las operaciones principales del negocio antes de deducir los gastos operativos.'
$data[71,2] = $null
$data[71,3] = 'number'
$data[71,4] = 'float'
$data[71,5] = 'Yes'
$data[71,6] = 'Yes'
$data[71,7] = $null
$data[71,8] = $null
$data[72,0] = 'income_statement.operating_expenses'
$data[72,1] = 'Los gastos totales incurridos por la empresa en sus actividades operativas normales, incluidos los gastos de venta, generales y administrativos.'
$data[72,2] = '122345.01'
$data[72,3] = 'number'
$data[72,4] = 'float'
$data[72,5] = 'Yes'
$data[72,6] = 'Yes'
$data[72,7] = $null
$data[72,8] = $null
$data[73,0] = 'income_statement.operating_income'
$data[73,1] = 'El beneficio obtenido de las operaciones principales del negocio después de deducir los gastos operativos, pero antes de considerar intereses, impuestos y otros elementos no operativos.'
$data[73,2] = '68545.01'
$data[73,3] = 'number'
$data[73,4] = 'float'
$data[73,5] = 'Yes'
$data[73,6] = 'Yes'
$data[73,7] = $null
$data[73,8] = $null
$data[74,0] = 'income_statement.operating_loss'
$data[74,1] = 'La pérdida incurrida por las operaciones principales del negocio después de deducir los gastos operativos, pero antes de considerar intereses, impuestos y otros elementos no operativos.'
$data[74,2] = $null
$data[74,3] = 'number'
$data[74,4] = 'float'
$data[74,5] = 'Yes'
$data[74,6] = 'Yes'
$data[74,7] = $null
$data[74,8] = $null
$data[75,0] = 'income_statement.financial_result'
$data[75,1] = 'El resultado neto de las actividades financieras, incluyendo los ingresos por intereses, los gastos por intereses y otras ganancias o pérdidas financieras.'
$data[75,2] = '15098.01'
$data[75,3] = 'number'
$data[75,4] = 'float'
$data[75,5] = 'Yes'
$data[75,6] = 'Yes'
$data[75,7] = $null
$data[75,8] = $null
$data[76,0] = 'income_statement.income_statement_financial_gains'
$data[76,1] = 'El ingreso financiero total positivo, incluyendo ingresos por intereses, ganancias por diferencias de cambio y otras ganancias de actividades de financiamiento. Este valor siempre debe ser positivo.'
$data[76,2] = '85000.01'
$data[76,3] = 'number'
$data[76,4] = 'float'
$data[76,5] = 'Yes'
$data[76,6] = 'Yes'
$data[76,7] = $null
$data[76,8] = $null
$data[77,0] = 'income_statement.income_statement_financial_costs'
$data[77,1] = 'Los gastos financieros totales, incluidos los gastos por intereses, las pérdidas por diferencias de cambio y otros costos incurridos por actividades de financiación. Este valor siempre debe ser negativo.'
$data[77,2] = '-32000.01'
$data[77,3] = 'number'
$data[77,4] = 'float'
$data[77,5] = 'Yes'
$data[77,6] = 'Yes'
$data[77,7] = $null
$data[77,8] = $null
$data[78,0] = 'income_statement.equity_in_earnings_of_affiliates'
$data[78,1] = 'La participación de la empresa en las ganancias o pérdidas de sus asociadas, entidades sobre las cuales tiene una influencia significativa pero no control.'
$data[78,2] = '5678.01'
$data[78,3] = 'number'
$data[78,4] = 'float'
$data[78,5] = 'Yes'
$data[78,6] = 'Yes'
$data[78,7] = $null
$data[78,8] = $null
$data[79,0] = 'income_statement.income_before_taxes'
$data[79,1] = 'El beneficio obtenido antes de contabilizar los gastos por impuesto sobre la renta.'
$data[79,2] = '89321.01'
$data[79,3] = 'number'
$data[79,4] = 'float'
$data[79,5] = 'Yes'
$data[79,6] = 'Yes'
$data[79,7] = $null
$data[79,8] = $null
$data[80,0] = 'income_statement.loss_before_taxes'
$data[80,1] = 'La pérdida incurrida antes de contabilizar los gastos por impuesto sobre la renta.'
$data[80,2] = $null
$data[80,3] = 'number'
$data[80,4] = 'float'
$data[80,5] = 'Yes'
$data[80,6] = 'Yes'
$data[80,7] = $null
$data[80,8] = $null
$data[81,0] = 'income_statement.income_taxes'
$data[81,1] = 'El monto total de los gastos por impuesto sobre la renta incurridos durante el período de informe.'
$data[81,2] = '20123.01'
$data[81,3] = 'number'
$data[81,4] = 'float'
$data[81,5] = 'Yes'
$data[81,6] = 'Yes'
$data[81,7] = $null
$data[81,8] = $null
$data[82,0] = 'income_statement.income_from_continuing_operations'
$data[82,1] = 'El beneficio obtenido de las operaciones comerciales en curso de la empresa después de deducir los gastos operativos e impuestos.'
$data[82,2] = '69198.01'
$data[82,3] = 'number'
$data[82,4] = 'float'
$data[82,5] = 'Yes'
$data[82,6] = 'Yes'
$data[82,7] = $null
$data[82,8] = $null
$data[83,0] = 'income_statement.loss_from_continuing_operations'
$data[83,1] = 'La pérdida incurrida por las operaciones comerciales en curso de la empresa después de deducir los gastos operativos e impuestos.'
$data[83,2] = $null
$data[83,3] = 'number'
$data[83,4] = 'float'
$data[83,5] = 'Yes'
$data[83,6] = 'Yes'
$data[83,7] = $null
$data[83,8] = $null
$data[84,0] = 'income_statement.discontinued_operations'
$data[84,1] = 'El resultado neto de las operaciones que han sido descontinuadas o vendidas durante el período de reporte.'
$data[84,2] = '0.01'
$data[84,3] = 'number'
$data[84,4] = 'float'
$data[84,5] = 'Yes'
$data[84,6] = 'Yes'
$data[84,7] = $null
$data[84,8] = $null
$data[85,0] = 'income_statement.net_income'
$data[85,1] = 'El beneficio total obtenido por la empresa después de deducir todos los gastos, incluidos los operativos, no operativos, intereses e impuestos.'
$data[85,2] = '69198.01'
$data[85,3] = 'number'
$data[85,4] = 'float'
$data[85,5] = 'Yes'
$data[85,6] = 'Yes'
$data[85,7] = $null
$data[85,8] = $null
$data[86,0] = 'income_statement.net_loss'
$data[86,1] = 'La pérdida total incurrida por la empresa después de deducir todos los gastos, incluidos los operativos, no operativos, intereses e impuestos.'
$data[86,2] = $null
$data[86,3] = 'number'
$data[86,4] = 'float'
$data[86,5] = 'Yes'
$data[86,6] = 'Yes'
$data[86,7] = $null
$data[86,8] = $null

# Set column C format to Text to preserve numeric-looking example strings exactly.
$ws.Range("C2:C88").NumberFormat = "@"

$ws.Range("A2:I88").Value = $data

# Step 3: Update the dimension / sanity check not required - Excel manages automatically.

# Step 4: Update text in the "links" sheet (credentials_storage row).
$wsLinks = $wb.Worksheets.Item("links")
$newText = 'Indica si se deben almacenar las credenciales (y la duración durante la cual se almacenarán las credenciales).

- Para enlaces recurrentes, esto se establece en `store` por defecto (y no se puede cambiar).
- Para enlaces únicos, esto se establece en `365d` por defecto.

Puede ser:
  - `store` para almacenar credenciales (hasta que se elimine el enlace)
  - `nostore` para no almacenar credenciales
  - Cualquier valor entre `1d` y `365d` para indicar el número de días que deseas que se almacenen las credenciales.

Para más información, consulta la sección <a href="https://developers.belvo.com/docs/data-retention-controls#credentials_storage" target="_blank">credentials_storage</a> de nuestro artículo sobre controles de retención de datos.'
$wsLinks.Cells.Item(12, 2).Value = $newText

Write-Host "Edit complete."